$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.841.02'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.09%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.40'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.20'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3616'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.74%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.06'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3242'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.129'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.77%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07020'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.994'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.91%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.41'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -5.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.661.89'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.560'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001039'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -5.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06600'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.36%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '78.77'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -4.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.878'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.78%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.62'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -7.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.55'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.830.06'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.36%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.27%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.399'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -11.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '147.26'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.49'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -6.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.845.33'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.69%  '

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.54'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.39%  '

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.188'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.94%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.100'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.697'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -12.44%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08440'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.72%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.660'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.17'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -9.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.282'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.85%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.120'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.11%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02240'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.67%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06001'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -7.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.247'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -5.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2048'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5884'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -6.05%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.767'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.65'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -5.15%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5573'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.39'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.67%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.930'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -5.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06965'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.188'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.54%  '
